# Weekly update: insert a new price observation for
# "Feria Lagunitas de Puerto Montt - Coliflor" as the new first data row
# (row 166), pushing the existing rows 166:181 down to 167:182.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(166).Insert()

$ws.Range("A166").Value = 4
$ws.Range("B166").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C166").Value = "Los Lagos"
$ws.Range("D166").Value = 44449
$ws.Range("E166").Value = 10
$ws.Range("F166").Value = 100112008
$ws.Range("G166").Value = "Coliflor"
$ws.Range("H166").Value = "Sin especificar"
$ws.Range("I166").Value = "Primera"
$ws.Range("J166").Value = 1000
$ws.Range("K166").Value = 1300
$ws.Range("L166").Value = 1300
$ws.Range("M166").Value = 1300
$ws.Range("N166").Value = "$/unidad"
$ws.Range("O166").Value = "Región Metropolitana"
$ws.Range("P166").Value = 1300
$ws.Range("Q166").Value = 1
$ws.Range("R166").Value = "Hortaliza"
